# Apply the point-by-point reply revisions described by the commit.
$d = $word.ActiveDocument

# 1) "streamlined both the abstract and the final paragraphs of the introduction"
#    -> "streamlined the final paragraph of the introduction"
$d.Content.Find.Execute(
    "streamlined both the abstract and the final paragraphs of the introduction",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "streamlined the final paragraph of the introduction", 2) | Out-Null

# 2) Remove the sentence about renaming 'Conlusions' -> merging paragraphs,
#    which the authors decided to cut entirely.
$openQuote = [char]8216
$closeQuote = [char]8217
$conclusionsSentence = "To improve the coherence of these sections, the " + $openQuote + "Conlusions" + $closeQuote + " consist of the last three paragraphs of our original discussion (instead of the last paragraph only). "
$d.Content.Find.Execute(
    $conclusionsSentence,
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 2) | Out-Null

# 3) "0.01). Therefore" -> "1%). Therefore"
$d.Content.Find.Execute(
    "0.01). Therefore",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "1%). Therefore", 2) | Out-Null

# 4) "We have revised the caption of Figure 2 to clarify ..."
#    -> "We have revised the last paragraph of our Methods section to clarify ..."
$d.Content.Find.Execute(
    "the caption of Figure 2 to clarify",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "the last paragraph of our Methods section to clarify", 2) | Out-Null
